$d = $word.ActiveDocument

# The placeholder/bookmark paragraph is the first paragraph in the body.
$p1 = $d.Paragraphs(1)

# Give the paragraph a (5 twip) border on all sides and bump its left
# indent from 120 -> 225 twips (6pt -> 11.25pt).
$p1.Range.ParagraphFormat.LeftIndent = 11.25
$borders = $p1.Range.ParagraphFormat.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# Replace the old placeholder id (plus its trailing space run) with the
# new placeholder id, collapsing the paragraph back down to a single run.
$p1.Range.Find.Execute("**ID__AFFARS_pgi_5315_topic_21__ID** ", $true, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_AFMC_PGI_5315_3D__ID**", 2)
